$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.075165666666667
$ws.Range("H2").Value = 9.225497000000001
$ws.Range("I2").Value = 0.02641273658732285
$ws.Range("J2").Value = 0.02641273658732285
$ws.Range("M2").Value = 0.668273
$ws.Range("N2").Value = 2.004819
$ws.Range("O2").Value = 0.01328414746766746
$ws.Range("P2").Value = 0.01328414746766746
$ws.Range("Q2").Value = 2.055050185560333
$ws.Range("R2").Value = 18.495451670043
$ws.Range("S2").Value = 0.0003508706878506524
$ws.Range("T2").Value = 0.0003508706878506523
$ws.Range("G3").Value = 3.075165666666667
$ws.Range("H3").Value = 9.225497000000001
$ws.Range("I3").Value = 0.02641273658732285
$ws.Range("J3").Value = 0.02641273658732285
$ws.Range("O3").Value = 0.3831531055114357
$ws.Range("P3").Value = 0.3831531055114357
$ws.Range("Q3").Value = 59.27357118669145
$ws.Range("R3").Value = 533.462140680223
$ws.Range("S3").Value = 0.01012012204848827
$ws.Range("T3").Value = 0.01012012204848827
$ws.Range("G4").Value = 3.075165666666667
$ws.Range("H4").Value = 9.225497000000001
$ws.Range("I4").Value = 0.02641273658732285
$ws.Range("J4").Value = 0.02641273658732285
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.6035627470208969
$ws.Range("P4").Value = 0.6035627470208967
$ws.Range("Q4").Value = 93.37081948853057
$ws.Range("R4").Value = 840.3373753967751
$ws.Range("S4").Value = 0.01594174385098392
$ws.Range("T4").Value = 0.01594174385098392
$ws.Range("I5").Value = 0.549422396165273
$ws.Range("J5").Value = 0.5494223961652731
$ws.Range("M5").Value = 0.668273
$ws.Range("N5").Value = 2.004819
$ws.Range("O5").Value = 0.01328414746766746
$ws.Range("P5").Value = 0.01328414746766746
$ws.Range("Q5").Value = 42.747959623857
$ws.Range("R5").Value = 384.731636614713
$ws.Range("S5").Value = 0.007298608132698699
$ws.Range("T5").Value = 0.007298608132698699
$ws.Range("I6").Value = 0.549422396165273
$ws.Range("J6").Value = 0.5494223961652731
$ws.Range("O6").Value = 0.3831531055114357
$ws.Range("P6").Value = 0.3831531055114357
$ws.Range("S6").Value = 0.2105128973282587
$ws.Range("T6").Value = 0.2105128973282587
$ws.Range("I7").Value = 0.549422396165273
$ws.Range("J7").Value = 0.5494223961652731
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.6035627470208969
$ws.Range("P7").Value = 0.6035627470208967
$ws.Range("R7").Value = 17480.20970304653
$ws.Range("S7").Value = 0.3316108907043157
$ws.Range("T7").Value = 0.3316108907043157
$ws.Range("I8").Value = 0.424164867247404
$ws.Range("J8").Value = 0.4241648672474041
$ws.Range("M8").Value = 0.668273
$ws.Range("N8").Value = 2.004819
$ws.Range("O8").Value = 0.01328414746766746
$ws.Range("P8").Value = 0.01328414746766746
$ws.Range("Q8").Value = 33.00226336877667
$ws.Range("R8").Value = 297.02037031899
$ws.Range("S8").Value = 0.005634668647118107
$ws.Range("T8").Value = 0.005634668647118106
$ws.Range("I9").Value = 0.424164867247404
$ws.Range("J9").Value = 0.4241648672474041
$ws.Range("O9").Value = 0.3831531055114357
$ws.Range("P9").Value = 0.3831531055114357
$ws.Range("S9").Value = 0.1625200861346887
$ws.Range("T9").Value = 0.1625200861346887
$ws.Range("I10").Value = 0.424164867247404
$ws.Range("J10").Value = 0.4241648672474041
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.6035627470208969
$ws.Range("P10").Value = 0.6035627470208967
$ws.Range("Q10").Value = 1499.451642286195
$ws.Range("S10").Value = 0.2560101124655972
$ws.Range("T10").Value = 0.2560101124655972
